$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "53.50", "0.000240").
# Force the cells to Text format before assigning so Excel preserves the
# literal digits/trailing zeros instead of coercing them into numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.765.59'
$ws.Range('E2').Value = '  -1.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.513.93'
$ws.Range('E3').Value = '  -1.69%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '614.94'
$ws.Range('E5').Value = '  +5.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '192.31'
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').Value = '  +0.71%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -3.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.664'
$ws.Range('E10').Value = '  +1.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.50'
$ws.Range('E11').Value = '  -2.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000309'
$ws.Range('E12').Value = '  -3.21%  '
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.081.42'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '618.81'
$ws.Range('E15').Value = '  +8.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.820.44'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('E18').Value = '  -1.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.522.16'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').Value = '  -1.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '109.42'
$ws.Range('E22').Value = '  +15.91%  '
$ws.Range('E23').Value = '  -3.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.72'
$ws.Range('E24').Value = '  +2.28%  '
$ws.Range('E25').Value = '  +3.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.10'
$ws.Range('E26').Value = '  +5.19%  '
$ws.Range('E27').Value = '  -1.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.72'
$ws.Range('E28').Value = '  +4.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.30'
$ws.Range('E29').Value = '  +5.04%  '
$ws.Range('E30').Value = '  -3.31%  '
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.94'
$ws.Range('E32').Value = '  +4.08%  '
$ws.Range('E33').Value = '  -0.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.58'
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('E35').Value = '  -5.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.674.31'
$ws.Range('E36').Value = '  +0.65%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('E38').Value = '  +6.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '516.60'
$ws.Range('E39').Value = '  -3.39%  '
$ws.Range('E40').Value = '  -4.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0781'
$ws.Range('E41').Value = '  -2.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.71'
$ws.Range('E42').Value = '  -4.70%  '
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('E46').Value = '  +3.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.33'
$ws.Range('E47').Value = '  -3.51%  '
$ws.Range('E48').Value = '  -5.78%  '
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.11'
$ws.Range('E50').Value = '  -1.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000240'
$ws.Range('E51').Value = '  -5.00%  '
